$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''22.436.49'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '''1.566.76'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').Value = '''284.69'
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('D7').Value = '''0.3621'
$ws.Range('E7').Value = '  -2.70%  '
$ws.Range('D8').Value = '''48.34'
$ws.Range('E8').Value = '  -3.10%  '
$ws.Range('D9').Value = '''0.3312'
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('D10').Value = '''1.119'
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('D11').Value = '''0.07371'
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').Value = '''1.002'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '''20.67'
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('D14').Value = '''5.937'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').Value = '''6.881'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '''1.579.74'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '''0.00001099'
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('D18').Value = '''87.47'
$ws.Range('E18').Value = '  -4.05%  '
$ws.Range('D19').Value = '''0.06708'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '''6.330'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = '''16.16'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '''11.97'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('D24').Value = '''22.438.36'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '''2.367'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').Value = '''2.527'
$ws.Range('E26').Value = '  -6.16%  '
$ws.Range('D27').Value = '''150.63'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = '''19.37'
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').Value = '''4.995'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '''123.84'
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('D31').Value = '''1.740.81'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').Value = '''1.026'
$ws.Range('E32').Value = '  -2.94%  '
$ws.Range('D33').Value = '''2.002'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').Value = '''6.070'
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').Value = '''9.657'
$ws.Range('E35').Value = '  -2.01%  '
$ws.Range('D36').Value = '''0.08232'
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('D37').Value = '''0.02398'
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('D38').Value = '''0.06397'
$ws.Range('E38').Value = '  -1.85%  '
$ws.Range('D39').Value = '''0.2223'
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''1.286'
$ws.Range('E40').Value = '  -3.54%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').Value = '''5.340'
$ws.Range('E41').Value = '  -2.30%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '''11.14'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '''0.6211'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''13.76'
$ws.Range('E44').Value = '  -1.51%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.5992'
$ws.Range('E45').Value = '  +3.01%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '''3.741'
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '''2.018'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''123.07'
$ws.Range('E48').Value = '  -5.67%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '''1.207'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.07191'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''75.34'
$ws.Range('E51').Value = '  -1.97%  '
